$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column ("Price") values are written as text. Coinranking prices are
# formatted strings (thousand separators as dots, fixed decimal places) and
# must stay t="inlineStr"/shared-string cells, not numeric cells. Assigning a
# numeric-looking string via .Value lets Excel auto-coerce it to a number
# (losing formatting, e.g. "592.01" -> 592.0099999999999238...). Forcing the
# cell to Text format before the write keeps it a string; resetting the
# style back to "Normal" afterwards avoids leaving a stray number-format
# style behind (matches the original, unstyled cells in the diff).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "68.254.86"
$ws.Range("E2").Value = "  +0.49%  "
Set-TextValue $ws.Range("D3") "3.735.00"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  +0.13%  "
Set-TextValue $ws.Range("D5") "592.01"
$ws.Range("E5").Value = "  -0.16%  "
Set-TextValue $ws.Range("D6") "166.07"
$ws.Range("E6").Value = "  +0.39%  "
Set-TextValue $ws.Range("D7") "3.733.33"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.61%  "
Set-TextValue $ws.Range("D10") "0.159"
$ws.Range("E10").Value = "  -0.21%  "
Set-TextValue $ws.Range("D11") "6.42"
$ws.Range("E11").Value = "  -0.20%  "
Set-TextValue $ws.Range("D12") "0.448"
$ws.Range("E12").Value = "  -0.05%  "
Set-TextValue $ws.Range("D13") "0.0000258"
$ws.Range("E13").Value = "  -2.23%  "
Set-TextValue $ws.Range("D14") "36.11"
$ws.Range("E14").Value = "  +0.55%  "
Set-TextValue $ws.Range("D15") "4.367.85"
$ws.Range("E15").Value = "  -0.15%  "
Set-TextValue $ws.Range("D16") "3.748.31"
$ws.Range("E16").Value = "  +0.25%  "
Set-TextValue $ws.Range("D17") "68.327.01"
$ws.Range("E17").Value = "  +0.76%  "
Set-TextValue $ws.Range("D18") "17.85"
$ws.Range("E18").Value = "  -2.76%  "
Set-TextValue $ws.Range("D19") "6.99"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  +0.43%  "
Set-TextValue $ws.Range("D22") "464.96"
$ws.Range("E22").Value = "  +0.06%  "
Set-TextValue $ws.Range("D23") "0.696"
$ws.Range("E23").Value = "  -0.90%  "
Set-TextValue $ws.Range("D24") "83.85"
$ws.Range("E24").Value = "  +1.13%  "
Set-TextValue $ws.Range("D25") "0.0000145"
$ws.Range("E25").Value = "  +7.19%  "
Set-TextValue $ws.Range("D26") "2.17"
$ws.Range("E26").Value = "  -0.67%  "
Set-TextValue $ws.Range("D27") "11.85"
$ws.Range("E27").Value = "  -0.94%  "
Set-TextValue $ws.Range("D28") "10.03"
$ws.Range("E28").Value = "  -1.46%  "
Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  -0.16%  "
Set-TextValue $ws.Range("D30") "3.895.46"
$ws.Range("E30").Value = "  +0.21%  "
Set-TextValue $ws.Range("D31") "2.76"
$ws.Range("E31").Value = "  -3.80%  "
Set-TextValue $ws.Range("D32") "7.27"
$ws.Range("E32").Value = "  -1.27%  "
Set-TextValue $ws.Range("D33") "29.72"
$ws.Range("E33").Value = "  -0.54%  "
Set-TextValue $ws.Range("D34") "2.15"
$ws.Range("E34").Value = "  -1.92%  "
Set-TextValue $ws.Range("D35") "9.14"
$ws.Range("E35").Value = "  +1.16%  "
Set-TextValue $ws.Range("D36") "0.999"
Set-TextValue $ws.Range("D37") "3.691.26"
$ws.Range("E37").Value = "  -0.12%  "
Set-TextValue $ws.Range("D38") "0.100"
$ws.Range("E38").Value = "  -1.61%  "
Set-TextValue $ws.Range("D39") "3.37"
$ws.Range("E39").Value = "  -2.67%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("E41").Value = "  +0.04%  "
Set-TextValue $ws.Range("D42") "5.76"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  -0.01%  "
Set-TextValue $ws.Range("D45") "43.75"
$ws.Range("E45").Value = "  +14.23%  "
Set-TextValue $ws.Range("D46") "0.300"
$ws.Range("E46").Value = "  -1.99%  "
Set-TextValue $ws.Range("D47") "46.57"
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("E48").Value = "  -0.25%  "
Set-TextValue $ws.Range("D49") "8.45"
$ws.Range("E49").Value = "  -1.02%  "
Set-TextValue $ws.Range("D50") "389.58"
$ws.Range("E50").Value = "  -1.32%  "
Set-TextValue $ws.Range("D51") "144.31"
$ws.Range("E51").Value = "  -0.27%  "
